$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) from existing header H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-15
$data = @(
    @(8,9),
    @(7,7),
    @(5,5),
    @(7,7),
    @(5,5),
    @(7,7),
    @(6,7),
    @(5,5),
    @(8,8),
    @(7,7),
    @(6,6),
    @(8,8),
    @(9,9),
    @(6,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
